$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 549-550 (existing rows 549+ shift down by 2)
$ws.Rows("549:550").Insert()

# New row 549
$ws.Range("A549").Value = 10
$ws.Range("B549").Value = "Vega Modelo de Temuco"
$ws.Range("C549").Value = "La Araucanía"
$ws.Range("D549").Value = 45106
$ws.Range("E549").Value = 9
$ws.Range("F549").Value = "Fruta"
$ws.Range("G549").Value = 100108
$ws.Range("H549").Value = "Tropicales y subtropicales"
$ws.Range("I549").Value = 100108002
$ws.Range("J549").Value = "Mango"
$ws.Range("K549").Value = "Sin especificar"
$ws.Range("L549").Value = "Primera"
$ws.Range("M549").Value = 900
$ws.Range("N549").Value = 10000
$ws.Range("O549").Value = 10000
$ws.Range("P549").Value = 10000
$ws.Range("Q549").Value = '$/bandeja 4 kilos'
$ws.Range("R549").Value = "Brasil"
$ws.Range("S549").Value = 2500
$ws.Range("T549").Value = 4

# New row 550
$ws.Range("A550").Value = 10
$ws.Range("B550").Value = "Vega Modelo de Temuco"
$ws.Range("C550").Value = "La Araucanía"
$ws.Range("D550").Value = 45106
$ws.Range("E550").Value = 9
$ws.Range("F550").Value = "Fruta"
$ws.Range("G550").Value = 100108
$ws.Range("H550").Value = "Tropicales y subtropicales"
$ws.Range("I550").Value = 100108002
$ws.Range("J550").Value = "Mango"
$ws.Range("K550").Value = "Sin especificar"
$ws.Range("L550").Value = "Primera"
$ws.Range("M550").Value = 1200
$ws.Range("N550").Value = 10000
$ws.Range("O550").Value = 10000
$ws.Range("P550").Value = 10000
$ws.Range("Q550").Value = '$/bandeja 4 kilos'
$ws.Range("R550").Value = "Perú"
$ws.Range("S550").Value = 2500
$ws.Range("T550").Value = 4

# Ensure date-format style (style index 2) applies to the new D-column cells,
# matching the rest of the Fecha column.
$ws.Range("D549").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D550").NumberFormat = "YYYY-MM-DD HH:MM:SS"
